$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Background Jobs" feature group
$ws.Range("B20").Value = "Background Jobs"
$ws.Range("C21").Value = "Excel sheet of montly expense to be processed by background job"
$ws.Range("C22").Value = "Upload file control"

# New "Backend Apis" feature group
$ws.Range("B24").Value = "Backend Apis"
$ws.Range("C24").Value = "Api Controllers need to be moved to a separate project"

# Widen column D to fit the new long text
$ws.Columns.Item(4).ColumnWidth = 25.3984375

# Update selection to mimic where the user left off editing
$ws.Range("C25").Select()
